$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 225.375
$ws.Range("I9").Value = 256.2
$ws.Range("K9").Value = 256.2
$ws.Range("M9").Value = -87.19999999999999

$ws.Range("H95").Value = 43641.5
$ws.Range("J95").Value = 43641.5
$ws.Range("L95").Value = 43641.5
$ws.Range("N95").Value = -49133.5

$ws.Range("H100").Value = 2369.2632
$ws.Range("I100").Value = 1780.7273
$ws.Range("J100").Value = 3178.5
$ws.Range("K100").Value = 1780.7273
$ws.Range("L100").Value = 3178.5
$ws.Range("M100").Value = -1239.7273
$ws.Range("N100").Value = -4260.5

$ws.Range("H137").Value = 2494.923
$ws.Range("I137").Value = 3045.4546
$ws.Range("J137").Value = 2091.2
$ws.Range("K137").Value = 9136.363799999999
$ws.Range("L137").Value = 6273.599999999999
$ws.Range("M137").Value = -6586.363799999999
$ws.Range("N137").Value = -11373.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2850.9546
$ws.Range("I2").Value = 1248.5834
$ws.Range("K2").Value = 1248.5834
$ws.Range("M2").Value = -1135.5834

$ws.Range("H28").Value = 14995
$ws.Range("I28").Value = 2012.6666
$ws.Range("K28").Value = 2012.6666
$ws.Range("M28").Value = -1820.6666

$ws.Range("H32").Value = 2054320
$ws.Range("I32").Value = 2363017.8
$ws.Range("K32").Value = 2363017.8
$ws.Range("M32").Value = -2362730.8

$ws.Range("H45").Value = 6614.7646
$ws.Range("I45").Value = 1863
$ws.Range("J45").Value = 9941
$ws.Range("K45").Value = 1863
$ws.Range("L45").Value = 9941
$ws.Range("M45").Value = -1486
$ws.Range("N45").Value = -10695

$ws.Range("H74").Value = 18028.77
$ws.Range("I74").Value = 44577.145
$ws.Range("J74").Value = 3161.68
$ws.Range("K74").Value = 44577.145
$ws.Range("L74").Value = 3161.68
$ws.Range("M74").Value = -43703.145
$ws.Range("N74").Value = -4909.68

$ws.Range("H77").Value = 18028.77
$ws.Range("I77").Value = 44577.145
$ws.Range("J77").Value = 3161.68
$ws.Range("K77").Value = 222885.725
$ws.Range("L77").Value = 15808.4
$ws.Range("M77").Value = -218517.725
$ws.Range("N77").Value = -24544.4

$ws.Range("H99").Value = 14995
$ws.Range("I99").Value = 2012.6666
$ws.Range("K99").Value = 2012.6666
$ws.Range("M99").Value = 982.3334

$ws.Range("H102").Value = 1815.5714
$ws.Range("I102").Value = 1942
$ws.Range("K102").Value = 1942
$ws.Range("M102").Value = -320

$ws.Range("H116").Value = 2850.9546
$ws.Range("I116").Value = 1248.5834
$ws.Range("K116").Value = 1248.5834
$ws.Range("M116").Value = 1045.4166

$ws.Range("H132").Value = 4479.5312
$ws.Range("I132").Value = 3326.36
$ws.Range("J132").Value = 8598
$ws.Range("K132").Value = 9979.08
$ws.Range("L132").Value = 25794
$ws.Range("M132").Value = -7449.08
$ws.Range("N132").Value = -30854

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2850.9546
$ws.Range("I3").Value = 1248.5834
$ws.Range("K3").Value = 1248.5834
$ws.Range("M3").Value = -1134.5834

$ws.Range("H107").Value = 45003560
$ws.Range("I107").Value = 59212900
$ws.Range("J107").Value = 7319.5
$ws.Range("K107").Value = 59212900
$ws.Range("L107").Value = 7319.5
$ws.Range("M107").Value = -59210980
$ws.Range("N107").Value = -11159.5

$ws.Range("H134").Value = 5430.804
$ws.Range("I134").Value = 2140.92
$ws.Range("K134").Value = 6422.76
$ws.Range("M134").Value = -3887.76

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9711.383
$ws.Range("I31").Value = 3810.6365
$ws.Range("J31").Value = 12533.479
$ws.Range("K31").Value = 3810.6365
$ws.Range("L31").Value = 12533.479
$ws.Range("M31").Value = -3515.6365
$ws.Range("N31").Value = -13123.479

$ws.Range("H34").Value = 9711.383
$ws.Range("I34").Value = 3810.6365
$ws.Range("J34").Value = 12533.479
$ws.Range("K34").Value = 3810.6365
$ws.Range("L34").Value = 12533.479
$ws.Range("M34").Value = -3608.6365
$ws.Range("N34").Value = -12937.479

$ws.Range("H86").Value = 12653700
$ws.Range("I86").Value = 20836332
$ws.Range("J86").Value = 379750
$ws.Range("K86").Value = 20836332
$ws.Range("L86").Value = 379750
$ws.Range("M86").Value = -20835209
$ws.Range("N86").Value = -381996

$ws.Range("H89").Value = 12653700
$ws.Range("I89").Value = 20836332
$ws.Range("J89").Value = 379750
$ws.Range("K89").Value = 104181660
$ws.Range("L89").Value = 1898750
$ws.Range("M89").Value = -104176044
$ws.Range("N89").Value = -1909982

$ws.Range("H99").Value = 8008.2
$ws.Range("I99").Value = 7967.6665
$ws.Range("K99").Value = 7967.6665
$ws.Range("M99").Value = -6469.6665

$ws.Range("H105").Value = 6498042.5
$ws.Range("I105").Value = 8930058
$ws.Range("J105").Value = 12666.667
$ws.Range("K105").Value = 8930058
$ws.Range("L105").Value = 12666.667
$ws.Range("M105").Value = -8928311
$ws.Range("N105").Value = -16160.667

$ws.Range("H122").Value = 1229.8462
$ws.Range("I122").Value = 962.5714
$ws.Range("J122").Value = 1541.6666
$ws.Range("K122").Value = 2887.7142
$ws.Range("L122").Value = 4624.9998
$ws.Range("M122").Value = -437.7142000000003
$ws.Range("N122").Value = -9524.9998

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H126").Value = 8008.2
$ws.Range("I126").Value = 7967.6665
$ws.Range("K126").Value = 23902.9995
$ws.Range("M126").Value = -21432.9995

$ws.Range("H134").Value = 6784.1
$ws.Range("I134").Value = 1691.7858
$ws.Range("K134").Value = 5075.357400000001
$ws.Range("M134").Value = -2540.357400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 990
$ws.Range("I8").Value = 990
$ws.Range("K8").Value = 2970
$ws.Range("M8").Value = -2831

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 128710.75
$ws.Range("I80").Value = 3666.3333
$ws.Range("K80").Value = 3666.3333
$ws.Range("M80").Value = -2668.3333

$ws.Range("H83").Value = 128710.75
$ws.Range("I83").Value = 3666.3333
$ws.Range("K83").Value = 18331.6665
$ws.Range("M83").Value = -13339.6665

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H132").Value = 6138.95
$ws.Range("I132").Value = 2037.3334
$ws.Range("J132").Value = 12291.375
$ws.Range("K132").Value = 6112.0002
$ws.Range("L132").Value = 36874.125
$ws.Range("M132").Value = -3582.0002
$ws.Range("N132").Value = -41934.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2009.2
$ws.Range("J46").Value = 3026
$ws.Range("L46").Value = 3026
$ws.Range("N46").Value = -3402

$ws.Range("H93").Value = 6286
$ws.Range("I93").Value = 5977.6
$ws.Range("J93").Value = 6671.5
$ws.Range("K93").Value = 5977.6
$ws.Range("L93").Value = 6671.5
$ws.Range("M93").Value = -4729.6
$ws.Range("N93").Value = -9167.5

$ws.Range("H100").Value = 4754.3335
$ws.Range("I100").Value = 2957.6
$ws.Range("K100").Value = 2957.6
$ws.Range("M100").Value = -2416.6

$ws.Range("H107").Value = 2087.125
$ws.Range("I107").Value = 2087.125
$ws.Range("K107").Value = 2087.125
$ws.Range("M107").Value = -167.125

$ws.Range("H122").Value = 4667.793
$ws.Range("I122").Value = 2687.2222
$ws.Range("K122").Value = 8061.6666
$ws.Range("M122").Value = -5611.6666

$ws.Range("H135").Value = 73759.28999999999
$ws.Range("J135").Value = 73759.28999999999
$ws.Range("L135").Value = 73759.28999999999
$ws.Range("N135").Value = -83899.28999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H107").Value = 27778522
$ws.Range("I107").Value = 446.5
$ws.Range("K107").Value = 1339.5
$ws.Range("M107").Value = 580.5

$ws.Range("I136").Value = 125001390
$ws.Range("K136").Value = 375004170
$ws.Range("M136").Value = -375001620

